$wb = $excel.ActiveWorkbook

# Update the parameter values on 'Sheet2'
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A2").Value = 12
$ws2.Range("B2").Value = 16
$ws2.Range("C2").Value = 18
$ws2.Range("D2").Value = 24

# Make 'Sheet2' the active sheet and select/activate cell E14 on it
$ws2.Activate()
$ws2.Range("E14").Select()

$wb.Save()
